# Add a new "LoginData" worksheet at the end of the workbook (after RegData),
# carrying old username/password test data (per commit message:
# "Old Selenium scripts added").

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing sheet so the final order is
# Sheet1, RegData, LoginData.
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "LoginData"

# Header row
$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "password"

# Data rows (old test credentials)
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "admin"

$ws.Range("A3").Value = "Admin123"
$ws.Range("B3").Value = "admin11"

$ws.Range("A4").Value = "Admin"
$ws.Range("B4").Value = "admin123"

# Highlight the header row, matching the gray-40% fill used on the header
$ws.Range("A1").Interior.ColorIndex = 55
$ws.Range("B1").Interior.ColorIndex = 55

# Size the header columns to fit their contents
$ws.Columns.Item(1).ColumnWidth = 10.2265625
$ws.Columns.Item(2).ColumnWidth = 9.62109375

# Keep RegData as the active/selected sheet, as it was before this edit.
$wb.Worksheets.Item("RegData").Activate()
